# Automatic update of files.
# The underlying source re-sorted/re-ordered two species observations that
# share the same location ("Sagviken, Dlr"); as a result the data
# previously on row 19 is now on row 20 and vice-versa. Swap the
# row-specific fields between row 19 and row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values of both rows ---
$row19 = @{
    A = $ws.Range("A19").Value2
    B = $ws.Range("B19").Value2
    D = $ws.Range("D19").Value2
    E = $ws.Range("E19").Value2
    F = $ws.Range("F19").Value2
    G = $ws.Range("G19").Value2
    H = $ws.Range("H19").Value2
    I = $ws.Range("I19").Value2
    J = $ws.Range("J19").Value2
    Q = $ws.Range("Q19").Value2
    R = $ws.Range("R19").Value2
    Z = $ws.Range("Z19").Value2
    AB = $ws.Range("AB19").Value2
}

$row20 = @{
    A = $ws.Range("A20").Value2
    B = $ws.Range("B20").Value2
    D = $ws.Range("D20").Value2
    E = $ws.Range("E20").Value2
    F = $ws.Range("F20").Value2
    G = $ws.Range("G20").Value2
    H = $ws.Range("H20").Value2
    I = $ws.Range("I20").Value2
    J = $ws.Range("J20").Value2
    Q = $ws.Range("Q20").Value2
    R = $ws.Range("R20").Value2
    Z = $ws.Range("Z20").Value2
    AB = $ws.Range("AB20").Value2
}

# Column I ("Antal") is stored as text in this sheet, even though its
# contents look numeric ("25", "20", ...). Excel's COM Value setter
# auto-detects numeric-looking strings and stores them as real numbers
# unless the cell is already Text-formatted, so flip the format to Text
# before assigning, then strip the formatting back off again so we don't
# leave a stray style on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- write row20's former values into row19 ---
$ws.Range("A19").Value = $row20.A
$ws.Range("B19").Value = $row20.B
$ws.Range("D19").Value = $row20.D
$ws.Range("E19").Value = $row20.E
$ws.Range("F19").Value = $row20.F
$ws.Range("G19").Value = $row20.G
$ws.Range("H19").Value = $row20.H
Set-TextValue $ws.Range("I19") $row20.I
$ws.Range("J19").Value = $row20.J
$ws.Range("Q19").Value = $row20.Q
$ws.Range("R19").Value = $row20.R
$ws.Range("Z19").Value = $row20.Z
$ws.Range("AB19").Value = $row20.AB

# --- write row19's former values into row20 ---
$ws.Range("A20").Value = $row19.A
$ws.Range("B20").Value = $row19.B
$ws.Range("D20").Value = $row19.D
$ws.Range("E20").Value = $row19.E
$ws.Range("F20").Value = $row19.F
$ws.Range("G20").Value = $row19.G
$ws.Range("H20").Value = $row19.H
Set-TextValue $ws.Range("I20") $row19.I
$ws.Range("J20").Value = $row19.J
$ws.Range("Q20").Value = $row19.Q
$ws.Range("R20").Value = $row19.R
$ws.Range("Z20").Value = $row19.Z
$ws.Range("AB20").Value = $row19.AB
